$d = $word.ActiveDocument

# 1) Hyperlink runs: add explicit "auto" font color (w:color w:val="auto").
#    wdColorAutomatic = -16777216. These are simple formatting-only edits
#    (no text mutation), so they do not disturb any other run in the
#    document.
$autoColor = -16777216
$d.Hyperlinks.Item(1).Range.Font.Color = $autoColor   # sofia.dutta17@gmail.com
$d.Hyperlinks.Item(2).Range.Font.Color = $autoColor   # https://linkedin.com/in/sofiadutta
$d.Hyperlinks.Item(3).Range.Font.Color = $autoColor   # https://sofiadutta.github.io

# 2) Insert the word "over " right before "a decade" in the summary
#    paragraph, turning "Software Engineer with a decade of experience..."
#    into "Software Engineer with over a decade of experience...".
#
#    A plain Range.InsertBefore() would cause the whole paragraph's runs
#    (which all already share identical visible formatting) to be
#    recombined into a single run. Wrapping the insertion as a tracked
#    change and then accepting just that one revision keeps each
#    surrounding run intact as its own <w:r> (matching how the paragraph
#    was originally split), while still leaving plain (non-revision) XML
#    behind once accepted.
$summaryParagraph = $d.Paragraphs.Item(6)
$summaryText = $summaryParagraph.Range.Text
$insertOffset = $summaryText.IndexOf("a decade")

$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true
try {
    $insertPos = $summaryParagraph.Range.Start + $insertOffset
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.InsertBefore("over ")
} finally {
    $d.TrackRevisions = $wasTracking
}

if ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}
